$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Fix the mislabeled "2050" column header (it was accidentally left holding
# a stray numeric value, 677.032366796629) on every table that has a 2050
# column, and drop the "Total" row that used to close out each table.
# ---------------------------------------------------------------------------

# Sheets 1-3: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
# "Atendimento a Ponta(MW)" -> header row is 2015 / 2030 / 2040 / 2050,
# Total row is row 13.
foreach ($idx in 1..3) {
    $ws = $wb.Worksheets.Item($idx)
    # E1 must read as text "2050" (not the stray number), matching the other
    # header cells which are also stored as text. A leading apostrophe forces
    # Excel to keep the numeric-looking literal as text.
    $ws.Range("E1").Value = "'2050"
    # Remove the old "Total" summary row entirely (shifts dimension A1:E13 -> A1:E12).
    $ws.Rows.Item(13).Delete()
}

# Sheet 4: "Potencia Incremental - SIN(MW)" -> header row is 2015 / 2015-2030 /
# 2031-2040 / 2041-2050, Total row is row 13.
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("E1").Value = "2041-2050"
$ws4.Rows.Item(13).Delete()

# Sheet 5: "Emissoes Totais (MtCO2eq)" -> header row is 2015 / 2030 / 2040 / 2050.
# This table never had a Total row, so only the header label needs fixing.
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("E1").Value = "'2050"

# Sheet 6: "Custo Total (bilhões de R$)" -> no 2050 column here (only A/B),
# just remove the trailing "Total" row (row 4), shrinking dimension A1:B4 -> A1:B3.
$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows.Item(4).Delete()
